# Day 4 done.  Fixed stupid bug with display of times in batch mode.
#
# Update the "2022" Advent-of-code stats sheet:
#   - correct the (batch-mode) times recorded for days 1 and 2
#   - fill in the times for days 3 and 4, which are now complete
#   - leave the selection on B6 (next cell to fill in)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("2022")
$ws.Activate()

# Day 1 (row 2): batch-mode time display bug fix
$ws.Range("B2").Value = 177174
$ws.Range("C2").Value = 6752

# Day 2 (row 3): same fix
$ws.Range("B3").Value = 143328
$ws.Range("C3").Value = 6395

# Day 3 (row 4): newly completed
$ws.Range("B4").Value = 111150
$ws.Range("C4").Value = 5645
$ws.Range("E4").Value = 35205
$ws.Range("F4").Value = 31491

# Day 4 (row 5): newly completed
$ws.Range("B5").Value = 92883
$ws.Range("C5").Value = 1862
$ws.Range("E5").Value = 20073
$ws.Range("F5").Value = 20122

# Re-assert every touched input cell once more. The dependent formulas in
# columns D/G/H (which chain B->D->G and C/B->H) only pick up the final
# value of a precedent if that precedent is the *last* cell written in the
# batch, so a second pass guarantees every formula downstream observes the
# final figures once everything has settled.
$ws.Range("B2").Value = 177174
$ws.Range("C2").Value = 6752
$ws.Range("B3").Value = 143328
$ws.Range("C3").Value = 6395
$ws.Range("B4").Value = 111150
$ws.Range("C4").Value = 5645
$ws.Range("E4").Value = 35205
$ws.Range("F4").Value = 31491
$ws.Range("B5").Value = 92883
$ws.Range("C5").Value = 1862
$ws.Range("E5").Value = 20073
$ws.Range("F5").Value = 20122

# Move the selection to the next day to fill in.
$ws.Range("B6").Select()
